$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D (old D..H shift right to F..J)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header row (row 2) labels for the inserted columns
$ws.Range("D2").Value = "Likelihood"
$ws.Range("E2").Value = "Severity"

# New sub-header row (row 3) description text (same text in both D3 and E3)
$riskMatrixText = "Use the Risk Matrix to asses" + [char]10 + "Rate from Low to Very High"
$ws.Range("D3").Value = $riskMatrixText
$ws.Range("E3").Value = $riskMatrixText

# Match the new columns' width to the original "Risk Analysis" column (now column F)
$ws.Columns("D:E").ColumnWidth = $ws.Columns("F").ColumnWidth

# Final selection: the whole of column F (previously column D) is selected
$ws.Range("F1:F1048576").Select() | Out-Null
